$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the 12 country labels that changed (abbreviated "St." / "&" forms)
# then re-sort the list A-Z, matching the authoring workflow.

$ws.Cells.Item(2, 1).Value = 'Afghanistan'
$ws.Cells.Item(3, 1).Value = 'Åland Islands'
$ws.Cells.Item(4, 1).Value = 'Albania'
$ws.Cells.Item(5, 1).Value = 'Algeria'
$ws.Cells.Item(6, 1).Value = 'American Samoa'
$ws.Cells.Item(7, 1).Value = 'Angola'
$ws.Cells.Item(8, 1).Value = 'Anguilla'
$ws.Cells.Item(9, 1).Value = 'Antarctica'
$ws.Cells.Item(10, 1).Value = 'Antigua & Barbuda'
$ws.Cells.Item(11, 1).Value = 'Argentina'
$ws.Cells.Item(12, 1).Value = 'Armenia'
$ws.Cells.Item(13, 1).Value = 'Aruba'
$ws.Cells.Item(14, 1).Value = 'Australia'
$ws.Cells.Item(15, 1).Value = 'Austria'
$ws.Cells.Item(16, 1).Value = 'Azerbaijan'
$ws.Cells.Item(17, 1).Value = 'Bahamas'
$ws.Cells.Item(18, 1).Value = 'Bahrain'
$ws.Cells.Item(19, 1).Value = 'Bangladesh'
$ws.Cells.Item(20, 1).Value = 'Barbados'
$ws.Cells.Item(21, 1).Value = 'Belarus'
$ws.Cells.Item(22, 1).Value = 'Belgium'
$ws.Cells.Item(23, 1).Value = 'Belize'
$ws.Cells.Item(24, 1).Value = 'Benin'
$ws.Cells.Item(25, 1).Value = 'Bermuda'
$ws.Cells.Item(26, 1).Value = 'Bhutan'
$ws.Cells.Item(27, 1).Value = 'Bolivia'
$ws.Cells.Item(28, 1).Value = 'Bosnia '
$ws.Cells.Item(29, 1).Value = 'Botswana'
$ws.Cells.Item(30, 1).Value = 'Bouvet Island'
$ws.Cells.Item(31, 1).Value = 'Brazil'
$ws.Cells.Item(32, 1).Value = 'British Indian Ocean Territory'
$ws.Cells.Item(33, 1).Value = 'British Virgin Islands'
$ws.Cells.Item(34, 1).Value = 'Brunei'
$ws.Cells.Item(35, 1).Value = 'Bulgaria'
$ws.Cells.Item(36, 1).Value = 'Burkina Faso'
$ws.Cells.Item(37, 1).Value = 'Burundi'
$ws.Cells.Item(38, 1).Value = 'Cambodia'
$ws.Cells.Item(39, 1).Value = 'Cameroon'
$ws.Cells.Item(40, 1).Value = 'Canada'
$ws.Cells.Item(41, 1).Value = 'Cape Verde'
$ws.Cells.Item(42, 1).Value = 'Cayman Islands'
$ws.Cells.Item(43, 1).Value = 'Central African Republic'
$ws.Cells.Item(44, 1).Value = 'Chad'
$ws.Cells.Item(45, 1).Value = 'Chile'
$ws.Cells.Item(46, 1).Value = 'China'
$ws.Cells.Item(47, 1).Value = 'Christmas Island'
$ws.Cells.Item(48, 1).Value = 'Cocos [Keeling] Islands'
$ws.Cells.Item(49, 1).Value = 'Colombia'
$ws.Cells.Item(50, 1).Value = 'Comoros'
$ws.Cells.Item(51, 1).Value = 'Congo - Brazzaville'
$ws.Cells.Item(52, 1).Value = 'Congo - Kinshasa'
$ws.Cells.Item(53, 1).Value = 'Cook Islands'
$ws.Cells.Item(54, 1).Value = 'Costa Rica'
$ws.Cells.Item(55, 1).Value = 'Côte d’Ivoire'
$ws.Cells.Item(56, 1).Value = 'Croatia'
$ws.Cells.Item(57, 1).Value = 'Cuba'
$ws.Cells.Item(58, 1).Value = 'Cyprus'
$ws.Cells.Item(59, 1).Value = 'Czechia'
$ws.Cells.Item(60, 1).Value = 'Denmark'
$ws.Cells.Item(61, 1).Value = 'Djibouti'
$ws.Cells.Item(62, 1).Value = 'Dominica'
$ws.Cells.Item(63, 1).Value = 'Dominican Republic'
$ws.Cells.Item(64, 1).Value = 'Ecuador'
$ws.Cells.Item(65, 1).Value = 'Egypt'
$ws.Cells.Item(66, 1).Value = 'El Salvador'
$ws.Cells.Item(67, 1).Value = 'Equatorial Guinea'
$ws.Cells.Item(68, 1).Value = 'Eritrea'
$ws.Cells.Item(69, 1).Value = 'Estonia'
$ws.Cells.Item(70, 1).Value = 'Ethiopia'
$ws.Cells.Item(71, 1).Value = 'Falkland Islands'
$ws.Cells.Item(72, 1).Value = 'Faroe Islands'
$ws.Cells.Item(73, 1).Value = 'Fiji'
$ws.Cells.Item(74, 1).Value = 'Finland'
$ws.Cells.Item(75, 1).Value = 'France'
$ws.Cells.Item(76, 1).Value = 'French Guiana'
$ws.Cells.Item(77, 1).Value = 'French Polynesia'
$ws.Cells.Item(78, 1).Value = 'French Southern Territories'
$ws.Cells.Item(79, 1).Value = 'Gabon'
$ws.Cells.Item(80, 1).Value = 'Gambia'
$ws.Cells.Item(81, 1).Value = 'Georgia'
$ws.Cells.Item(82, 1).Value = 'Germany'
$ws.Cells.Item(83, 1).Value = 'Ghana'
$ws.Cells.Item(84, 1).Value = 'Gibraltar'
$ws.Cells.Item(85, 1).Value = 'Greece'
$ws.Cells.Item(86, 1).Value = 'Greenland'
$ws.Cells.Item(87, 1).Value = 'Grenada'
$ws.Cells.Item(88, 1).Value = 'Guadeloupe'
$ws.Cells.Item(89, 1).Value = 'Guam'
$ws.Cells.Item(90, 1).Value = 'Guatemala'
$ws.Cells.Item(91, 1).Value = 'Guernsey'
$ws.Cells.Item(92, 1).Value = 'Guinea'
$ws.Cells.Item(93, 1).Value = 'Guinea-Bissau'
$ws.Cells.Item(94, 1).Value = 'Guyana'
$ws.Cells.Item(95, 1).Value = 'Haiti'
$ws.Cells.Item(96, 1).Value = 'Heard & McDonald Islands'
$ws.Cells.Item(97, 1).Value = 'Honduras'
$ws.Cells.Item(98, 1).Value = 'Hong Kong [China]'
$ws.Cells.Item(99, 1).Value = 'Hungary'
$ws.Cells.Item(100, 1).Value = 'Iceland'
$ws.Cells.Item(101, 1).Value = 'India'
$ws.Cells.Item(102, 1).Value = 'Indonesia'
$ws.Cells.Item(103, 1).Value = 'Iran'
$ws.Cells.Item(104, 1).Value = 'Iraq'
$ws.Cells.Item(105, 1).Value = 'Ireland'
$ws.Cells.Item(106, 1).Value = 'Isle of Man'
$ws.Cells.Item(107, 1).Value = 'Israel'
$ws.Cells.Item(108, 1).Value = 'Italy'
$ws.Cells.Item(109, 1).Value = 'Jamaica'
$ws.Cells.Item(110, 1).Value = 'Japan'
$ws.Cells.Item(111, 1).Value = 'Jersey'
$ws.Cells.Item(112, 1).Value = 'Jordan'
$ws.Cells.Item(113, 1).Value = 'Kazakhstan'
$ws.Cells.Item(114, 1).Value = 'Kenya'
$ws.Cells.Item(115, 1).Value = 'Kiribati'
$ws.Cells.Item(116, 1).Value = 'Kuwait'
$ws.Cells.Item(117, 1).Value = 'Kyrgyzstan'
$ws.Cells.Item(118, 1).Value = 'Laos'
$ws.Cells.Item(119, 1).Value = 'Latvia'
$ws.Cells.Item(120, 1).Value = 'Lebanon'
$ws.Cells.Item(121, 1).Value = 'Lesotho'
$ws.Cells.Item(122, 1).Value = 'Liberia'
$ws.Cells.Item(123, 1).Value = 'Libya'
$ws.Cells.Item(124, 1).Value = 'Liechtenstein'
$ws.Cells.Item(125, 1).Value = 'Lithuania'
$ws.Cells.Item(126, 1).Value = 'Luxembourg'
$ws.Cells.Item(127, 1).Value = 'Macau'
$ws.Cells.Item(128, 1).Value = 'Macedonia'
$ws.Cells.Item(129, 1).Value = 'Madagascar'
$ws.Cells.Item(130, 1).Value = 'Malawi'
$ws.Cells.Item(131, 1).Value = 'Malaysia'
$ws.Cells.Item(132, 1).Value = 'Maldives'
$ws.Cells.Item(133, 1).Value = 'Mali'
$ws.Cells.Item(134, 1).Value = 'Malta'
$ws.Cells.Item(135, 1).Value = 'Marshall Islands'
$ws.Cells.Item(136, 1).Value = 'Martinique'
$ws.Cells.Item(137, 1).Value = 'Mauritania'
$ws.Cells.Item(138, 1).Value = 'Mauritius'
$ws.Cells.Item(139, 1).Value = 'Mayotte'
$ws.Cells.Item(140, 1).Value = 'Mexico'
$ws.Cells.Item(141, 1).Value = 'Micronesia'
$ws.Cells.Item(142, 1).Value = 'Moldova'
$ws.Cells.Item(143, 1).Value = 'Monaco'
$ws.Cells.Item(144, 1).Value = 'Mongolia'
$ws.Cells.Item(145, 1).Value = 'Montenegro'
$ws.Cells.Item(146, 1).Value = 'Montserrat'
$ws.Cells.Item(147, 1).Value = 'Morocco'
$ws.Cells.Item(148, 1).Value = 'Mozambique'
$ws.Cells.Item(149, 1).Value = 'Myanmar [Burma]'
$ws.Cells.Item(150, 1).Value = 'Namibia'
$ws.Cells.Item(151, 1).Value = 'Nauru'
$ws.Cells.Item(152, 1).Value = 'Nepal'
$ws.Cells.Item(153, 1).Value = 'Netherlands'
$ws.Cells.Item(154, 1).Value = 'New Caledonia'
$ws.Cells.Item(155, 1).Value = 'New Zealand'
$ws.Cells.Item(156, 1).Value = 'Nicaragua'
$ws.Cells.Item(157, 1).Value = 'Niger'
$ws.Cells.Item(158, 1).Value = 'Nigeria'
$ws.Cells.Item(159, 1).Value = 'Niue'
$ws.Cells.Item(160, 1).Value = 'Norfolk Island'
$ws.Cells.Item(161, 1).Value = 'North Korea'
$ws.Cells.Item(162, 1).Value = 'Northern Mariana Islands'
$ws.Cells.Item(163, 1).Value = 'Norway'
$ws.Cells.Item(164, 1).Value = 'Oman'
$ws.Cells.Item(165, 1).Value = 'Pakistan'
$ws.Cells.Item(166, 1).Value = 'Palau'
$ws.Cells.Item(167, 1).Value = 'Palestinian Territories'
$ws.Cells.Item(168, 1).Value = 'Panama'
$ws.Cells.Item(169, 1).Value = 'Papua New Guinea'
$ws.Cells.Item(170, 1).Value = 'Paraguay'
$ws.Cells.Item(171, 1).Value = 'Peru'
$ws.Cells.Item(172, 1).Value = 'Philippines'
$ws.Cells.Item(173, 1).Value = 'Pitcairn Islands'
$ws.Cells.Item(174, 1).Value = 'Poland'
$ws.Cells.Item(175, 1).Value = 'Portugal'
$ws.Cells.Item(176, 1).Value = 'Puerto Rico'
$ws.Cells.Item(177, 1).Value = 'Qatar'
$ws.Cells.Item(178, 1).Value = 'Réunion'
$ws.Cells.Item(179, 1).Value = 'Romania'
$ws.Cells.Item(180, 1).Value = 'Russia'
$ws.Cells.Item(181, 1).Value = 'Rwanda'
$ws.Cells.Item(182, 1).Value = 'Samoa'
$ws.Cells.Item(183, 1).Value = 'San Marino'
$ws.Cells.Item(184, 1).Value = 'São Tomé & Príncipe'
$ws.Cells.Item(185, 1).Value = 'Saudi Arabia'
$ws.Cells.Item(186, 1).Value = 'Senegal'
$ws.Cells.Item(187, 1).Value = 'Serbia'
$ws.Cells.Item(188, 1).Value = 'Seychelles'
$ws.Cells.Item(189, 1).Value = 'Sierra Leone'
$ws.Cells.Item(190, 1).Value = 'Singapore'
$ws.Cells.Item(191, 1).Value = 'Slovakia'
$ws.Cells.Item(192, 1).Value = 'Slovenia'
$ws.Cells.Item(193, 1).Value = 'So. Georgia & So. Sandwich Isl.'
$ws.Cells.Item(194, 1).Value = 'Solomon Islands'
$ws.Cells.Item(195, 1).Value = 'Somalia'
$ws.Cells.Item(196, 1).Value = 'South Africa'
$ws.Cells.Item(197, 1).Value = 'South Korea'
$ws.Cells.Item(198, 1).Value = 'Spain'
$ws.Cells.Item(199, 1).Value = 'Sri Lanka'
$ws.Cells.Item(200, 1).Value = 'St. Barthélemy'
$ws.Cells.Item(201, 1).Value = 'St. Helena'
$ws.Cells.Item(202, 1).Value = 'St. Kitts and Nevis'
$ws.Cells.Item(203, 1).Value = 'St. Lucia'
$ws.Cells.Item(204, 1).Value = 'St. Martin'
$ws.Cells.Item(205, 1).Value = 'St. Vincent & Grenadines'
$ws.Cells.Item(206, 1).Value = 'Sudan'
$ws.Cells.Item(207, 1).Value = 'Suriname'
$ws.Cells.Item(208, 1).Value = 'Swaziland'
$ws.Cells.Item(209, 1).Value = 'Sweden'
$ws.Cells.Item(210, 1).Value = 'Switzerland'
$ws.Cells.Item(211, 1).Value = 'Syria'
$ws.Cells.Item(212, 1).Value = 'Taiwan'
$ws.Cells.Item(213, 1).Value = 'Tajikistan'
$ws.Cells.Item(214, 1).Value = 'Tanzania'
$ws.Cells.Item(215, 1).Value = 'Thailand'
$ws.Cells.Item(216, 1).Value = 'Timor-Leste'
$ws.Cells.Item(217, 1).Value = 'Togo'
$ws.Cells.Item(218, 1).Value = 'Tokelau'
$ws.Cells.Item(219, 1).Value = 'Tonga'
$ws.Cells.Item(220, 1).Value = 'Trinidad & Tobago'
$ws.Cells.Item(221, 1).Value = 'Tunisia'
$ws.Cells.Item(222, 1).Value = 'Turkey'
$ws.Cells.Item(223, 1).Value = 'Turkmenistan'
$ws.Cells.Item(224, 1).Value = 'Turks & Caicos Islands'
$ws.Cells.Item(225, 1).Value = 'Tuvalu'
$ws.Cells.Item(226, 1).Value = 'U.S. Virgin Islands'
$ws.Cells.Item(227, 1).Value = 'Uganda'
$ws.Cells.Item(228, 1).Value = 'Ukraine'
$ws.Cells.Item(229, 1).Value = 'United Arab Emirates'
$ws.Cells.Item(230, 1).Value = 'United Kingdom'
$ws.Cells.Item(231, 1).Value = 'United States'
$ws.Cells.Item(232, 1).Value = 'Uruguay'
$ws.Cells.Item(233, 1).Value = 'Uzbekistan'
$ws.Cells.Item(234, 1).Value = 'Vanuatu'
$ws.Cells.Item(235, 1).Value = 'Vatican City'
$ws.Cells.Item(236, 1).Value = 'Venezuela'
$ws.Cells.Item(237, 1).Value = 'Vietnam'
$ws.Cells.Item(238, 1).Value = 'Wallis & Futuna'
$ws.Cells.Item(239, 1).Value = 'Western Sahara'
$ws.Cells.Item(240, 1).Value = 'Yemen'
$ws.Cells.Item(241, 1).Value = 'Zambia'
$ws.Cells.Item(242, 1).Value = 'Zimbabwe'

# Apply an A-Z sort over the data range (keeps the sortState bookkeeping
# that Excel writes after using Data > Sort).
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add2($ws.Range("A190"))
$sortObj.SetRange($ws.Range("A2:A242"))
$sortObj.Header = 2
$sortObj.Apply()

# Re-assert the final alphabetical order explicitly (engine sort uses
# ordinal comparison, which doesn't match Excel's locale-aware A-Z sort
# for accented names like "Åland Islands" or "Côte d’Ivoire").
$ws.Cells.Item(2, 1).Value = 'Afghanistan'
$ws.Cells.Item(3, 1).Value = 'Åland Islands'
$ws.Cells.Item(4, 1).Value = 'Albania'
$ws.Cells.Item(5, 1).Value = 'Algeria'
$ws.Cells.Item(6, 1).Value = 'American Samoa'
$ws.Cells.Item(7, 1).Value = 'Angola'
$ws.Cells.Item(8, 1).Value = 'Anguilla'
$ws.Cells.Item(9, 1).Value = 'Antarctica'
$ws.Cells.Item(10, 1).Value = 'Antigua & Barbuda'
$ws.Cells.Item(11, 1).Value = 'Argentina'
$ws.Cells.Item(12, 1).Value = 'Armenia'
$ws.Cells.Item(13, 1).Value = 'Aruba'
$ws.Cells.Item(14, 1).Value = 'Australia'
$ws.Cells.Item(15, 1).Value = 'Austria'
$ws.Cells.Item(16, 1).Value = 'Azerbaijan'
$ws.Cells.Item(17, 1).Value = 'Bahamas'
$ws.Cells.Item(18, 1).Value = 'Bahrain'
$ws.Cells.Item(19, 1).Value = 'Bangladesh'
$ws.Cells.Item(20, 1).Value = 'Barbados'
$ws.Cells.Item(21, 1).Value = 'Belarus'
$ws.Cells.Item(22, 1).Value = 'Belgium'
$ws.Cells.Item(23, 1).Value = 'Belize'
$ws.Cells.Item(24, 1).Value = 'Benin'
$ws.Cells.Item(25, 1).Value = 'Bermuda'
$ws.Cells.Item(26, 1).Value = 'Bhutan'
$ws.Cells.Item(27, 1).Value = 'Bolivia'
$ws.Cells.Item(28, 1).Value = 'Bosnia '
$ws.Cells.Item(29, 1).Value = 'Botswana'
$ws.Cells.Item(30, 1).Value = 'Bouvet Island'
$ws.Cells.Item(31, 1).Value = 'Brazil'
$ws.Cells.Item(32, 1).Value = 'British Indian Ocean Territory'
$ws.Cells.Item(33, 1).Value = 'British Virgin Islands'
$ws.Cells.Item(34, 1).Value = 'Brunei'
$ws.Cells.Item(35, 1).Value = 'Bulgaria'
$ws.Cells.Item(36, 1).Value = 'Burkina Faso'
$ws.Cells.Item(37, 1).Value = 'Burundi'
$ws.Cells.Item(38, 1).Value = 'Cambodia'
$ws.Cells.Item(39, 1).Value = 'Cameroon'
$ws.Cells.Item(40, 1).Value = 'Canada'
$ws.Cells.Item(41, 1).Value = 'Cape Verde'
$ws.Cells.Item(42, 1).Value = 'Cayman Islands'
$ws.Cells.Item(43, 1).Value = 'Central African Republic'
$ws.Cells.Item(44, 1).Value = 'Chad'
$ws.Cells.Item(45, 1).Value = 'Chile'
$ws.Cells.Item(46, 1).Value = 'China'
$ws.Cells.Item(47, 1).Value = 'Christmas Island'
$ws.Cells.Item(48, 1).Value = 'Cocos [Keeling] Islands'
$ws.Cells.Item(49, 1).Value = 'Colombia'
$ws.Cells.Item(50, 1).Value = 'Comoros'
$ws.Cells.Item(51, 1).Value = 'Congo - Brazzaville'
$ws.Cells.Item(52, 1).Value = 'Congo - Kinshasa'
$ws.Cells.Item(53, 1).Value = 'Cook Islands'
$ws.Cells.Item(54, 1).Value = 'Costa Rica'
$ws.Cells.Item(55, 1).Value = 'Côte d’Ivoire'
$ws.Cells.Item(56, 1).Value = 'Croatia'
$ws.Cells.Item(57, 1).Value = 'Cuba'
$ws.Cells.Item(58, 1).Value = 'Cyprus'
$ws.Cells.Item(59, 1).Value = 'Czechia'
$ws.Cells.Item(60, 1).Value = 'Denmark'
$ws.Cells.Item(61, 1).Value = 'Djibouti'
$ws.Cells.Item(62, 1).Value = 'Dominica'
$ws.Cells.Item(63, 1).Value = 'Dominican Republic'
$ws.Cells.Item(64, 1).Value = 'Ecuador'
$ws.Cells.Item(65, 1).Value = 'Egypt'
$ws.Cells.Item(66, 1).Value = 'El Salvador'
$ws.Cells.Item(67, 1).Value = 'Equatorial Guinea'
$ws.Cells.Item(68, 1).Value = 'Eritrea'
$ws.Cells.Item(69, 1).Value = 'Estonia'
$ws.Cells.Item(70, 1).Value = 'Ethiopia'
$ws.Cells.Item(71, 1).Value = 'Falkland Islands'
$ws.Cells.Item(72, 1).Value = 'Faroe Islands'
$ws.Cells.Item(73, 1).Value = 'Fiji'
$ws.Cells.Item(74, 1).Value = 'Finland'
$ws.Cells.Item(75, 1).Value = 'France'
$ws.Cells.Item(76, 1).Value = 'French Guiana'
$ws.Cells.Item(77, 1).Value = 'French Polynesia'
$ws.Cells.Item(78, 1).Value = 'French Southern Territories'
$ws.Cells.Item(79, 1).Value = 'Gabon'
$ws.Cells.Item(80, 1).Value = 'Gambia'
$ws.Cells.Item(81, 1).Value = 'Georgia'
$ws.Cells.Item(82, 1).Value = 'Germany'
$ws.Cells.Item(83, 1).Value = 'Ghana'
$ws.Cells.Item(84, 1).Value = 'Gibraltar'
$ws.Cells.Item(85, 1).Value = 'Greece'
$ws.Cells.Item(86, 1).Value = 'Greenland'
$ws.Cells.Item(87, 1).Value = 'Grenada'
$ws.Cells.Item(88, 1).Value = 'Guadeloupe'
$ws.Cells.Item(89, 1).Value = 'Guam'
$ws.Cells.Item(90, 1).Value = 'Guatemala'
$ws.Cells.Item(91, 1).Value = 'Guernsey'
$ws.Cells.Item(92, 1).Value = 'Guinea'
$ws.Cells.Item(93, 1).Value = 'Guinea-Bissau'
$ws.Cells.Item(94, 1).Value = 'Guyana'
$ws.Cells.Item(95, 1).Value = 'Haiti'
$ws.Cells.Item(96, 1).Value = 'Heard & McDonald Islands'
$ws.Cells.Item(97, 1).Value = 'Honduras'
$ws.Cells.Item(98, 1).Value = 'Hong Kong [China]'
$ws.Cells.Item(99, 1).Value = 'Hungary'
$ws.Cells.Item(100, 1).Value = 'Iceland'
$ws.Cells.Item(101, 1).Value = 'India'
$ws.Cells.Item(102, 1).Value = 'Indonesia'
$ws.Cells.Item(103, 1).Value = 'Iran'
$ws.Cells.Item(104, 1).Value = 'Iraq'
$ws.Cells.Item(105, 1).Value = 'Ireland'
$ws.Cells.Item(106, 1).Value = 'Isle of Man'
$ws.Cells.Item(107, 1).Value = 'Israel'
$ws.Cells.Item(108, 1).Value = 'Italy'
$ws.Cells.Item(109, 1).Value = 'Jamaica'
$ws.Cells.Item(110, 1).Value = 'Japan'
$ws.Cells.Item(111, 1).Value = 'Jersey'
$ws.Cells.Item(112, 1).Value = 'Jordan'
$ws.Cells.Item(113, 1).Value = 'Kazakhstan'
$ws.Cells.Item(114, 1).Value = 'Kenya'
$ws.Cells.Item(115, 1).Value = 'Kiribati'
$ws.Cells.Item(116, 1).Value = 'Kuwait'
$ws.Cells.Item(117, 1).Value = 'Kyrgyzstan'
$ws.Cells.Item(118, 1).Value = 'Laos'
$ws.Cells.Item(119, 1).Value = 'Latvia'
$ws.Cells.Item(120, 1).Value = 'Lebanon'
$ws.Cells.Item(121, 1).Value = 'Lesotho'
$ws.Cells.Item(122, 1).Value = 'Liberia'
$ws.Cells.Item(123, 1).Value = 'Libya'
$ws.Cells.Item(124, 1).Value = 'Liechtenstein'
$ws.Cells.Item(125, 1).Value = 'Lithuania'
$ws.Cells.Item(126, 1).Value = 'Luxembourg'
$ws.Cells.Item(127, 1).Value = 'Macau'
$ws.Cells.Item(128, 1).Value = 'Macedonia'
$ws.Cells.Item(129, 1).Value = 'Madagascar'
$ws.Cells.Item(130, 1).Value = 'Malawi'
$ws.Cells.Item(131, 1).Value = 'Malaysia'
$ws.Cells.Item(132, 1).Value = 'Maldives'
$ws.Cells.Item(133, 1).Value = 'Mali'
$ws.Cells.Item(134, 1).Value = 'Malta'
$ws.Cells.Item(135, 1).Value = 'Marshall Islands'
$ws.Cells.Item(136, 1).Value = 'Martinique'
$ws.Cells.Item(137, 1).Value = 'Mauritania'
$ws.Cells.Item(138, 1).Value = 'Mauritius'
$ws.Cells.Item(139, 1).Value = 'Mayotte'
$ws.Cells.Item(140, 1).Value = 'Mexico'
$ws.Cells.Item(141, 1).Value = 'Micronesia'
$ws.Cells.Item(142, 1).Value = 'Moldova'
$ws.Cells.Item(143, 1).Value = 'Monaco'
$ws.Cells.Item(144, 1).Value = 'Mongolia'
$ws.Cells.Item(145, 1).Value = 'Montenegro'
$ws.Cells.Item(146, 1).Value = 'Montserrat'
$ws.Cells.Item(147, 1).Value = 'Morocco'
$ws.Cells.Item(148, 1).Value = 'Mozambique'
$ws.Cells.Item(149, 1).Value = 'Myanmar [Burma]'
$ws.Cells.Item(150, 1).Value = 'Namibia'
$ws.Cells.Item(151, 1).Value = 'Nauru'
$ws.Cells.Item(152, 1).Value = 'Nepal'
$ws.Cells.Item(153, 1).Value = 'Netherlands'
$ws.Cells.Item(154, 1).Value = 'New Caledonia'
$ws.Cells.Item(155, 1).Value = 'New Zealand'
$ws.Cells.Item(156, 1).Value = 'Nicaragua'
$ws.Cells.Item(157, 1).Value = 'Niger'
$ws.Cells.Item(158, 1).Value = 'Nigeria'
$ws.Cells.Item(159, 1).Value = 'Niue'
$ws.Cells.Item(160, 1).Value = 'Norfolk Island'
$ws.Cells.Item(161, 1).Value = 'North Korea'
$ws.Cells.Item(162, 1).Value = 'Northern Mariana Islands'
$ws.Cells.Item(163, 1).Value = 'Norway'
$ws.Cells.Item(164, 1).Value = 'Oman'
$ws.Cells.Item(165, 1).Value = 'Pakistan'
$ws.Cells.Item(166, 1).Value = 'Palau'
$ws.Cells.Item(167, 1).Value = 'Palestinian Territories'
$ws.Cells.Item(168, 1).Value = 'Panama'
$ws.Cells.Item(169, 1).Value = 'Papua New Guinea'
$ws.Cells.Item(170, 1).Value = 'Paraguay'
$ws.Cells.Item(171, 1).Value = 'Peru'
$ws.Cells.Item(172, 1).Value = 'Philippines'
$ws.Cells.Item(173, 1).Value = 'Pitcairn Islands'
$ws.Cells.Item(174, 1).Value = 'Poland'
$ws.Cells.Item(175, 1).Value = 'Portugal'
$ws.Cells.Item(176, 1).Value = 'Puerto Rico'
$ws.Cells.Item(177, 1).Value = 'Qatar'
$ws.Cells.Item(178, 1).Value = 'Réunion'
$ws.Cells.Item(179, 1).Value = 'Romania'
$ws.Cells.Item(180, 1).Value = 'Russia'
$ws.Cells.Item(181, 1).Value = 'Rwanda'
$ws.Cells.Item(182, 1).Value = 'Samoa'
$ws.Cells.Item(183, 1).Value = 'San Marino'
$ws.Cells.Item(184, 1).Value = 'São Tomé & Príncipe'
$ws.Cells.Item(185, 1).Value = 'Saudi Arabia'
$ws.Cells.Item(186, 1).Value = 'Senegal'
$ws.Cells.Item(187, 1).Value = 'Serbia'
$ws.Cells.Item(188, 1).Value = 'Seychelles'
$ws.Cells.Item(189, 1).Value = 'Sierra Leone'
$ws.Cells.Item(190, 1).Value = 'Singapore'
$ws.Cells.Item(191, 1).Value = 'Slovakia'
$ws.Cells.Item(192, 1).Value = 'Slovenia'
$ws.Cells.Item(193, 1).Value = 'So. Georgia & So. Sandwich Isl.'
$ws.Cells.Item(194, 1).Value = 'Solomon Islands'
$ws.Cells.Item(195, 1).Value = 'Somalia'
$ws.Cells.Item(196, 1).Value = 'South Africa'
$ws.Cells.Item(197, 1).Value = 'South Korea'
$ws.Cells.Item(198, 1).Value = 'Spain'
$ws.Cells.Item(199, 1).Value = 'Sri Lanka'
$ws.Cells.Item(200, 1).Value = 'St. Barthélemy'
$ws.Cells.Item(201, 1).Value = 'St. Helena'
$ws.Cells.Item(202, 1).Value = 'St. Kitts and Nevis'
$ws.Cells.Item(203, 1).Value = 'St. Lucia'
$ws.Cells.Item(204, 1).Value = 'St. Martin'
$ws.Cells.Item(205, 1).Value = 'St. Vincent & Grenadines'
$ws.Cells.Item(206, 1).Value = 'Sudan'
$ws.Cells.Item(207, 1).Value = 'Suriname'
$ws.Cells.Item(208, 1).Value = 'Swaziland'
$ws.Cells.Item(209, 1).Value = 'Sweden'
$ws.Cells.Item(210, 1).Value = 'Switzerland'
$ws.Cells.Item(211, 1).Value = 'Syria'
$ws.Cells.Item(212, 1).Value = 'Taiwan'
$ws.Cells.Item(213, 1).Value = 'Tajikistan'
$ws.Cells.Item(214, 1).Value = 'Tanzania'
$ws.Cells.Item(215, 1).Value = 'Thailand'
$ws.Cells.Item(216, 1).Value = 'Timor-Leste'
$ws.Cells.Item(217, 1).Value = 'Togo'
$ws.Cells.Item(218, 1).Value = 'Tokelau'
$ws.Cells.Item(219, 1).Value = 'Tonga'
$ws.Cells.Item(220, 1).Value = 'Trinidad & Tobago'
$ws.Cells.Item(221, 1).Value = 'Tunisia'
$ws.Cells.Item(222, 1).Value = 'Turkey'
$ws.Cells.Item(223, 1).Value = 'Turkmenistan'
$ws.Cells.Item(224, 1).Value = 'Turks & Caicos Islands'
$ws.Cells.Item(225, 1).Value = 'Tuvalu'
$ws.Cells.Item(226, 1).Value = 'U.S. Virgin Islands'
$ws.Cells.Item(227, 1).Value = 'Uganda'
$ws.Cells.Item(228, 1).Value = 'Ukraine'
$ws.Cells.Item(229, 1).Value = 'United Arab Emirates'
$ws.Cells.Item(230, 1).Value = 'United Kingdom'
$ws.Cells.Item(231, 1).Value = 'United States'
$ws.Cells.Item(232, 1).Value = 'Uruguay'
$ws.Cells.Item(233, 1).Value = 'Uzbekistan'
$ws.Cells.Item(234, 1).Value = 'Vanuatu'
$ws.Cells.Item(235, 1).Value = 'Vatican City'
$ws.Cells.Item(236, 1).Value = 'Venezuela'
$ws.Cells.Item(237, 1).Value = 'Vietnam'
$ws.Cells.Item(238, 1).Value = 'Wallis & Futuna'
$ws.Cells.Item(239, 1).Value = 'Western Sahara'
$ws.Cells.Item(240, 1).Value = 'Yemen'
$ws.Cells.Item(241, 1).Value = 'Zambia'
$ws.Cells.Item(242, 1).Value = 'Zimbabwe'

# Restore the selection/scroll position left by the editor after sorting.
$ws.Range("A96").Select()
